# Add the new rows 4-10 to the daily cash-register sheet (caja-2023-09-05).
# All existing cells in the sheet are stored as TEXT (even the numeric
# looking ones such as prices / weights), so every new cell is forced to
# Text number-format before the value is written, to avoid Excel silently
# re-typing "20000" etc. as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("ESTAMPA",  "20000", "GATO",  "15", "CREDITO"),
    @("NUTRIBON", "8300",  "PERRO", "20", "TRANSFERENCIA"),
    @("DOGUI",    "8000",  "PERRO", "21", "TRANSFERENCIA"),
    @("NUTRIBON", "13333", "GATO",  "20", "EFECTIVO"),
    @("NUTRIBON", "5100",  "PERRO", "20", "EFECTIVO"),
    @("DOGUI",    "499",   "PERRO", "2",  "EFECTIVO"),
    @("JHOLA",    "5100",  "PERRO", "20", "EFECTIVO")
)

$startRow = 4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $rowRange = $ws.Range("A" + $r + ":E" + $r)
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
